# BoM update: refresh item names/costs, zero-out a cancelled line (red
# highlight), add two new BoM rows (LM7805, SlipRing), and push the
# Total row down to make room, widening its SUM ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Carry the "Total" row's formatting down to row 11 before row 8's
#    cells get reused for the new SlipRing line.
# ---------------------------------------------------------------------
$ws.Range("D8:F8").Copy()
$ws.Range("D11:F11").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) New item names (shared-string order follows this edit order).
# ---------------------------------------------------------------------
$ws.Range("B7").Value = "LM7805"
$ws.Range("B8").Value = "SlipRing"
$ws.Range("B2").Value = "Ultrasonic Sensor (HC-SR04)"
$ws.Range("B6").Value = "Logic Level Shifter (TXB0104PWR + BOB)"
$ws.Range("B5").Value = "IR sensor for Homing"

# ---------------------------------------------------------------------
# 3) Row 2: E2 was "=13"; it's now a plain literal value.
# ---------------------------------------------------------------------
$ws.Range("E2").Value = 13

# ---------------------------------------------------------------------
# 4) Row 4 (Magnetic Encoders) cancelled: zero cost/total, red highlight.
# ---------------------------------------------------------------------
$tmpl = $ws.Range("H1")
$tmpl.Interior.Color = 192
$tmpl.Font.Name = "Aptos Narrow"
$tmpl.Copy()
$ws.Range("A4:F4").PasteSpecial(-4122)
$tmpl.Clear()

$ws.Range("E4").Value = 0
$ws.Range("F4").Formula = "=C4/D4*E4*0"

# ---------------------------------------------------------------------
# 5) Row 5 (IR sensor for Homing): new qty-with-purchase / cost.
# ---------------------------------------------------------------------
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 13

# ---------------------------------------------------------------------
# 6) Row 6 (Logic Level Shifter): new cost.
# ---------------------------------------------------------------------
$ws.Range("E6").Value = 5.32

# ---------------------------------------------------------------------
# 7) Rows 7 & 8 are new BoM lines -- give them the same look as the
#    other data rows (A col style, B:F col style) before filling values.
# ---------------------------------------------------------------------
$ws.Range("A3:F3").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122)
$ws.Range("A3:F3").Copy()
$ws.Range("A8:F8").PasteSpecial(-4122)

$ws.Range("A7").Value = 7
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 2.82
$ws.Range("F7").Formula = "=C7/D7*E7"

$ws.Range("A8").Value = 8
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 15
$ws.Range("F8").Formula = "=C8/D8*E8"

# ---------------------------------------------------------------------
# 8) Total row, now at row 11, spanning the widened data range.
# ---------------------------------------------------------------------
$ws.Range("D11").Value = "Total"
$ws.Range("E11").Formula = "=SUM(E2:E8)"
$ws.Range("F11").Formula = "=SUM(F2:F8)"

# ---------------------------------------------------------------------
# 9) Sheet-view bookkeeping (matches the authored selection/zoom).
# ---------------------------------------------------------------------
$ws.Range("E15").Select()
$excel.ActiveWindow.Zoom = 137
